$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2212.6282
$ws.Range("I15").Value = 2212.6282
$ws.Range("K15").Value = 6637.8846
$ws.Range("M15").Value = -6468.8846
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null
$ws.Range("H111").Value = 7374
$ws.Range("J111").Value = 6493.8335
$ws.Range("L111").Value = 19481.5005
$ws.Range("N111").Value = -25615.5005
$ws.Range("H116").Value = 3076.45
$ws.Range("I116").Value = 2747.6365
$ws.Range("J116").Value = 3478.3333
$ws.Range("K116").Value = 2747.6365
$ws.Range("L116").Value = 3478.3333
$ws.Range("M116").Value = 694.3634999999999
$ws.Range("N116").Value = -10362.3333
$ws.Range("H129").Value = 1096.4615
$ws.Range("I129").Value = 495.45456
$ws.Range("J129").Value = 1257.7073
$ws.Range("K129").Value = 1486.36368
$ws.Range("L129").Value = 3773.1219
$ws.Range("M129").Value = 3513.63632
$ws.Range("N129").Value = -13773.1219
$ws.Range("H138").Value = 2143.3604
$ws.Range("J138").Value = 2021.3
$ws.Range("L138").Value = 6063.9
$ws.Range("N138").Value = -16343.9
$ws.Range("H141").Value = 6528.871
$ws.Range("I141").Value = 2461.389
$ws.Range("J141").Value = 12160.77
$ws.Range("K141").Value = 7384.167
$ws.Range("L141").Value = 36482.31
$ws.Range("M141").Value = -2204.167
$ws.Range("N141").Value = -46842.31

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 861032.3
$ws.Range("I32").Value = 1121555.8
$ws.Range("J32").Value = 21568.223
$ws.Range("K32").Value = 1121555.8
$ws.Range("L32").Value = 21568.223
$ws.Range("M32").Value = -1121268.8
$ws.Range("N32").Value = -22142.223
$ws.Range("H37").Value = 11796
$ws.Range("J37").Value = 11796
$ws.Range("L37").Value = 11796
$ws.Range("N37").Value = -12342
$ws.Range("H44").Value = 18098.375
$ws.Range("J44").Value = 18098.375
$ws.Range("L44").Value = 18098.375
$ws.Range("N44").Value = -19074.375
$ws.Range("H45").Value = 2370.4443
$ws.Range("I45").Value = 1598.1428
$ws.Range("J45").Value = 3202.1538
$ws.Range("K45").Value = 1598.1428
$ws.Range("L45").Value = 3202.1538
$ws.Range("M45").Value = -1221.1428
$ws.Range("N45").Value = -3956.1538
$ws.Range("H55").Value = 29999
$ws.Range("J55").Value = 29999
$ws.Range("L55").Value = 29999
$ws.Range("N55").Value = -30629
$ws.Range("H61").Value = 2334.8647
$ws.Range("I61").Value = 1943
$ws.Range("J61").Value = 4014.2856
$ws.Range("K61").Value = 1943
$ws.Range("L61").Value = 4014.2856
$ws.Range("M61").Value = -1731
$ws.Range("N61").Value = -4438.2856
$ws.Range("H80").Value = 21249
$ws.Range("J80").Value = 21249
$ws.Range("L80").Value = 21249
$ws.Range("N80").Value = -23245
$ws.Range("H83").Value = 21249
$ws.Range("J83").Value = 21249
$ws.Range("L83").Value = 63747
$ws.Range("N83").Value = -73731
$ws.Range("H132").Value = 5124.273
$ws.Range("I132").Value = 5721.4546
$ws.Range("J132").Value = 4527.091
$ws.Range("K132").Value = 17164.3638
$ws.Range("L132").Value = 13581.273
$ws.Range("M132").Value = -14634.3638
$ws.Range("N132").Value = -18641.273
$ws.Range("H136").Value = 2334.8647
$ws.Range("I136").Value = 1943
$ws.Range("J136").Value = 4014.2856
$ws.Range("K136").Value = 5829
$ws.Range("L136").Value = 12042.8568
$ws.Range("M136").Value = -3279
$ws.Range("N136").Value = -17142.8568

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = $null
$ws.Range("H4").Value = 26230.924
$ws.Range("J4").Value = 26230.924
$ws.Range("L4").Value = 26230.924
$ws.Range("N4").Value = -26454.924
$ws.Range("H31").Value = 6250.9
$ws.Range("I31").Value = 1414.5518
$ws.Range("K31").Value = 1414.5518
$ws.Range("M31").Value = -1119.5518
$ws.Range("H34").Value = 6250.9
$ws.Range("I34").Value = 1414.5518
$ws.Range("K34").Value = 1414.5518
$ws.Range("M34").Value = -1212.5518

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8752830
$ws.Range("I4").Value = 13336693
$ws.Range("J4").Value = 7695015.5
$ws.Range("K4").Value = 40010079
$ws.Range("L4").Value = 23085046.5
$ws.Range("M4").Value = -40009967
$ws.Range("N4").Value = -23085270.5
$ws.Range("H50").Value = 18518828
$ws.Range("I50").Value = 246.41667
$ws.Range("J50").Value = 55555990
$ws.Range("K50").Value = 739.25001
$ws.Range("L50").Value = 166667970
$ws.Range("M50").Value = -258.25001
$ws.Range("N50").Value = -166668932
$ws.Range("H53").Value = 18518828
$ws.Range("I53").Value = 246.41667
$ws.Range("J53").Value = 55555990
$ws.Range("K53").Value = 739.25001
$ws.Range("L53").Value = 166667970
$ws.Range("M53").Value = -258.25001
$ws.Range("N53").Value = -166668932
$ws.Range("H113").Value = 831.9706
$ws.Range("I113").Value = 467.7857
$ws.Range("J113").Value = 1086.9
$ws.Range("K113").Value = 1403.3571
$ws.Range("L113").Value = 3260.7
$ws.Range("M113").Value = 766.6428999999998
$ws.Range("N113").Value = -7600.700000000001
$ws.Range("H131").Value = 1041.909
$ws.Range("I131").Value = 310
$ws.Range("J131").Value = 1115.1
$ws.Range("K131").Value = 930
$ws.Range("L131").Value = 3345.3
$ws.Range("M131").Value = 4110
$ws.Range("N131").Value = -13425.3

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = $null
$ws.Range("H70").Value = 5365.041
$ws.Range("I70").Value = 5364.516
$ws.Range("J70").Value = 5365.9443
$ws.Range("K70").Value = 5364.516
$ws.Range("L70").Value = 5365.9443
$ws.Range("M70").Value = -5094.516
$ws.Range("N70").Value = -5905.9443
$ws.Range("H73").Value = 5365.041
$ws.Range("I73").Value = 5364.516
$ws.Range("J73").Value = 5365.9443
$ws.Range("K73").Value = 5364.516
$ws.Range("L73").Value = 5365.9443
$ws.Range("M73").Value = -4428.516
$ws.Range("N73").Value = -7237.9443
$ws.Range("H132").Value = 2738.0715
$ws.Range("I132").Value = 2416.0908
$ws.Range("J132").Value = 3918.6667
$ws.Range("K132").Value = 7248.2724
$ws.Range("L132").Value = 11756.0001
$ws.Range("M132").Value = -4718.2724
$ws.Range("N132").Value = -16816.0001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 28111.555
$ws.Range("J2").Value = 28111.555
$ws.Range("L2").Value = 28111.555
$ws.Range("N2").Value = -28335.555
$ws.Range("H68").Value = 5083.3335
$ws.Range("I68").Value = 4000
$ws.Range("J68").Value = 5300
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 5300
$ws.Range("M68").Value = -3251
$ws.Range("N68").Value = -6798
$ws.Range("H71").Value = 5083.3335
$ws.Range("I71").Value = 4000
$ws.Range("J71").Value = 5300
$ws.Range("K71").Value = 20000
$ws.Range("L71").Value = 26500
$ws.Range("M71").Value = -16256
$ws.Range("N71").Value = -33988
$ws.Range("H132").Value = 3646.1365
$ws.Range("I132").Value = 3044.3572
$ws.Range("J132").Value = 4699.25
$ws.Range("K132").Value = 9133.071599999999
$ws.Range("L132").Value = 14097.75
$ws.Range("M132").Value = -6603.071599999999
$ws.Range("N132").Value = -19157.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 95000
$ws.Range("J64").Value = 95000
$ws.Range("L64").Value = 95000
$ws.Range("N64").Value = -95496
$ws.Range("H67").Value = 95000
$ws.Range("J67").Value = 95000
$ws.Range("L67").Value = 95000
$ws.Range("N67").Value = -96716
$ws.Range("H107").Value = 799
$ws.Range("I107").Value = 778.8
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 2336.4
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = -416.3999999999996
$ws.Range("N107").Value = -6540
$ws.Range("H122").Value = 2961.525
$ws.Range("I122").Value = 1520.0435
$ws.Range("J122").Value = 4911.7646
$ws.Range("K122").Value = 4560.1305
$ws.Range("L122").Value = 14735.2938
$ws.Range("M122").Value = -2110.1305
$ws.Range("N122").Value = -19635.2938
$ws.Range("H132").Value = 9262959
$ws.Range("I132").Value = 10668.333
$ws.Range("J132").Value = 11113417
$ws.Range("K132").Value = 32004.999
$ws.Range("L132").Value = 33340251
$ws.Range("M132").Value = -29474.999
$ws.Range("N132").Value = -33345311
$ws.Range("H136").Value = 1942.8462
$ws.Range("I136").Value = 1649.25
$ws.Range("J136").Value = 3285
$ws.Range("K136").Value = 4947.75
$ws.Range("L136").Value = 9855
$ws.Range("M136").Value = -2397.75
$ws.Range("N136").Value = -14955
$ws.Range("H138").Value = 80422
$ws.Range("J138").Value = 80422
$ws.Range("L138").Value = 80422
$ws.Range("N138").Value = -90702
$ws.Range("H139").Value = 80705
$ws.Range("J139").Value = 80705
$ws.Range("L139").Value = 80705
$ws.Range("N139").Value = -90985
